# Update "想去人数" (number of people wanting to go) figures on the
# "展览" and "全部类型" sheets for rows 2-5 (F column).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 9665
    $ws.Range("F3").Value = 214
    $ws.Range("F4").Value = 33
    $ws.Range("F5").Value = 547
}
